$d = $word.ActiveDocument

# The paragraph originally reads (split across 3 runs, with a grammar
# proofing mark wrapping the middle run):
#   "- Teleconsultation Access → Grants " + "doctors" + " permission to conduct video consultations.  "
# Collapse it into a single contiguous run of text (Word's Find/Replace
# naturally merges the runs it rewrites and drops the now-irrelevant
# proofErr markers around "doctors").
$found = $d.Content.Find.Execute(
    "Grants doctors permission",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Grants doctors permission", 2)
